$d = $word.ActiveDocument

# The document ends with 4 empty trailing paragraphs right before the
# final section break. The 3rd of those four (index Count-1) is where
# the new "Sitio en netlify..." block gets typed in, followed by a
# series of brand-new paragraphs (text lines, blank spacer lines and
# hyperlinks) that were not present before. The very last trailing
# empty paragraph is left untouched at the end.

$targetIndex = $d.Paragraphs.Count - 1
$target = $d.Paragraphs.Item($targetIndex)

# 1) "Sitio en netlify (...)"
$target.Range.InsertAfter("Sitio en netlify (ir checando hasta cuando dura… parece que un mes)")

# helper: always add a new empty paragraph after the paragraph at $idx
# and return the index of the freshly created (still empty) paragraph.
function New-ParaAfter([int]$idx) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    return $idx + 1
}

# 2) blank line
$idx = New-ParaAfter $targetIndex

# 3) blank line before the first hyperlink
$idx = New-ParaAfter $idx

# 4) hyperlink: https://app.netlify.com/teams/mapaznavarro/sites
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$url1 = "https://app.netlify.com/teams/mapaznavarro/sites"
$p.Range.InsertAfter($url1)
$p = $d.Paragraphs.Item($idx)
$linkRange = $d.Range($p.Range.Start, $p.Range.End - 1)
$d.Hyperlinks.Add($linkRange, $url1) | Out-Null
$p = $d.Paragraphs.Item($idx)
$styleFix = $d.Range($p.Range.Start, $p.Range.End - 1)
$styleFix.Find.ClearFormatting()
$styleFix.Find.Replacement.ClearFormatting()
$styleFix.Find.Replacement.Style = "EnlacedeInternet"
$styleFix.Find.Execute($url1, $false, $false, $false, $false, $false, $true, 1, $false, $url1, 2) | Out-Null

# 5) blank line
$idx = New-ParaAfter $idx

# 6) "Ruta de publicacion:"
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Ruta de publicacion:")

# 7) hyperlink https://gorgeous-phoenix-45c3f9.netlify.app/ followed by a
#    trailing run containing a single space.
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$url2 = "https://gorgeous-phoenix-45c3f9.netlify.app/"
$p.Range.InsertAfter($url2)
$p = $d.Paragraphs.Item($idx)
$linkRange2 = $d.Range($p.Range.Start, $p.Range.End - 1)
$d.Hyperlinks.Add($linkRange2, $url2) | Out-Null
$p = $d.Paragraphs.Item($idx)
$styleFix2 = $d.Range($p.Range.Start, $p.Range.End - 1)
$styleFix2.Find.ClearFormatting()
$styleFix2.Find.Replacement.ClearFormatting()
$styleFix2.Find.Replacement.Style = "EnlacedeInternet"
$styleFix2.Find.Execute($url2, $false, $false, $false, $false, $false, $true, 1, $false, $url2, 2) | Out-Null
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter(" ")

# 8) blank line
$idx = New-ParaAfter $idx

# 9) "Ruta fuentes en Git"
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("Ruta fuentes en Git")

# 10) hyperlink https://github.com/Mapaznavarro/LlamaApiGeolocalizaIP
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$url3 = "https://github.com/Mapaznavarro/LlamaApiGeolocalizaIP"
$p.Range.InsertAfter($url3)
$p = $d.Paragraphs.Item($idx)
$linkRange3 = $d.Range($p.Range.Start, $p.Range.End - 1)
$d.Hyperlinks.Add($linkRange3, $url3) | Out-Null
$p = $d.Paragraphs.Item($idx)
$styleFix3 = $d.Range($p.Range.Start, $p.Range.End - 1)
$styleFix3.Find.ClearFormatting()
$styleFix3.Find.Replacement.ClearFormatting()
$styleFix3.Find.Replacement.Style = "EnlacedeInternet"
$styleFix3.Find.Execute($url3, $false, $false, $false, $false, $false, $true, 1, $false, $url3, 2) | Out-Null

# 11) blank line
$idx = New-ParaAfter $idx

# 12) "IMPORTANTE: ..."
$idx = New-ParaAfter $idx
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("IMPORTANTE: La rama main está conectada a la versión que corre en el sitio netlify.")

Write-Host "Done. Paragraph count now: $($d.Paragraphs.Count)"
